$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: change numeric values to inline strings
$ws.Range("A3").Value = "Pd"
$ws.Range("B3").Value = "Cd"
$ws.Range("C3").Value = "Ru"

# Row 4
$ws.Range("A4").Value = 1
$ws.Range("B4").Value = 11
$ws.Range("C4").Value = 21

# Row 5
$ws.Range("A5").Value = 2
$ws.Range("B5").Value = "#N/A"
$ws.Range("C5").Value = 22

# Row 6
$ws.Range("A6").Value = 3
$ws.Range("B6").Value = 13
$ws.Range("C6").Value = 23
